$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2375886524822695
$ws.Range("C2").Value = 0.4964539007092199
$ws.Range("J2").Value = 0.01773049645390071
$ws.Range("P2").Value = 0.1524822695035461
$ws.Range("S2").Value = 0.09574468085106383

# Row 3
$ws.Range("B3").Value = 0.006711409395973154
$ws.Range("C3").Value = 0.04697986577181208
$ws.Range("J3").Value = 0.02684563758389262
$ws.Range("P3").Value = 0.7785234899328859
$ws.Range("S3").Value = 0.1409395973154362

# Row 4
$ws.Range("P4").Value = 0.6949152542372882
$ws.Range("S4").Value = 0.3050847457627119

# Row 5
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25

# Row 6
$ws.Range("B6").Value = 0.09836065573770492
$ws.Range("D6").Value = 0.01639344262295082
$ws.Range("F6").Value = 0.03825136612021858
$ws.Range("J6").Value = 0.2513661202185792
$ws.Range("O6").Value = 0.02185792349726776
$ws.Range("Q6").Value = 0.1967213114754098
$ws.Range("R6").Value = 0.07103825136612021
$ws.Range("S6").Value = 0.3060109289617486

# Row 7
$ws.Range("B7").Value = 0.08173076923076923
$ws.Range("D7").Value = 0.04326923076923077
$ws.Range("F7").Value = 0.05288461538461538
$ws.Range("J7").Value = 0.1346153846153846
$ws.Range("O7").Value = 0.01923076923076923
$ws.Range("Q7").Value = 0.2019230769230769
$ws.Range("R7").Value = 0.1057692307692308
$ws.Range("S7").Value = 0.3605769230769231

# Row 8
$ws.Range("B8").Value = 0.07413010590015129
$ws.Range("D8").Value = 0.0226928895612708
$ws.Range("E8").Value = 0.00453857791225416
$ws.Range("F8").Value = 0.03630862329803328
$ws.Range("J8").Value = 0.113464447806354
$ws.Range("O8").Value = 0.02874432677760968
$ws.Range("Q8").Value = 0.1800302571860817
$ws.Range("R8").Value = 0.1119515885022693
$ws.Range("S8").Value = 0.4281391830559758

# Row 9
$ws.Range("B9").Value = 0.1142857142857143
$ws.Range("D9").Value = 0.005714285714285714
$ws.Range("F9").Value = 0.05142857142857143
$ws.Range("J9").Value = 0.1085714285714286
$ws.Range("O9").Value = 0.01714285714285714
$ws.Range("Q9").Value = 0.2057142857142857
$ws.Range("R9").Value = 0.07428571428571429
$ws.Range("S9").Value = 0.4228571428571429

# Row 10
$ws.Range("B10").Value = 0.08830950378469302
$ws.Range("D10").Value = 0.02775441547518924
$ws.Range("E10").Value = 0.0008410428931875525
$ws.Range("F10").Value = 0.06728343145500421
$ws.Range("J10").Value = 0.1068124474348192
$ws.Range("O10").Value = 0.0159798149705635
$ws.Range("Q10").Value = 0.2354920100925147
$ws.Range("R10").Value = 0.09503784693019345
$ws.Range("S10").Value = 0.3624894869638352

# Row 11
$ws.Range("G11").Value = 0.09621993127147767
$ws.Range("J11").Value = 0.08934707903780069
$ws.Range("K11").Value = 0.134020618556701
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("S11").Value = 0.01374570446735395

# Row 12
$ws.Range("G12").Value = 0.7638190954773869
$ws.Range("J12").Value = 0.1959798994974874
$ws.Range("K12").Value = 0.005025125628140704
$ws.Range("L12").Value = 0.01005025125628141
$ws.Range("S12").Value = 0.02512562814070352

# Row 13
$ws.Range("G13").Value = 0.74
$ws.Range("J13").Value = 0.18
$ws.Range("S13").Value = 0.08

# Row 15
$ws.Range("F15").Value = 0.004587155963302753
$ws.Range("H15").Value = 0.2293577981651376
$ws.Range("I15").Value = 0.05504587155963303
$ws.Range("J15").Value = 0.3165137614678899
$ws.Range("K15").Value = 0.04587155963302753
$ws.Range("M15").Value = 0.01834862385321101
$ws.Range("O15").Value = 0.05504587155963303
$ws.Range("S15").Value = 0.2752293577981652

# Row 16
$ws.Range("F16").Value = 0.005128205128205128
$ws.Range("H16").Value = 0.2205128205128205
$ws.Range("I16").Value = 0.09230769230769231
$ws.Range("J16").Value = 0.3692307692307693
$ws.Range("K16").Value = 0.1435897435897436
$ws.Range("M16").Value = 0.02564102564102564
$ws.Range("N16").Value = 0.01025641025641026
$ws.Range("O16").Value = 0.06153846153846154
$ws.Range("S16").Value = 0.07179487179487179

# Row 17
$ws.Range("F17").Value = 0.01565557729941291
$ws.Range("H17").Value = 0.2661448140900196
$ws.Range("I17").Value = 0.07240704500978473
$ws.Range("J17").Value = 0.3424657534246575
$ws.Range("K17").Value = 0.09197651663405088
$ws.Range("M17").Value = 0.02152641878669276
$ws.Range("O17").Value = 0.07436399217221135
$ws.Range("S17").Value = 0.1154598825831702

# Row 18
$ws.Range("F18").Value = 0.01716738197424893
$ws.Range("H18").Value = 0.2918454935622318
$ws.Range("I18").Value = 0.0815450643776824
$ws.Range("J18").Value = 0.3776824034334764
$ws.Range("K18").Value = 0.07725321888412018
$ws.Range("M18").Value = 0.0128755364806867
$ws.Range("O18").Value = 0.02575107296137339
$ws.Range("S18").Value = 0.1158798283261803

# Row 19
$ws.Range("F19").Value = 0.007142857142857143
$ws.Range("H19").Value = 0.2928571428571429
$ws.Range("I19").Value = 0.07222222222222222
$ws.Range("J19").Value = 0.3301587301587302
$ws.Range("K19").Value = 0.1150793650793651
$ws.Range("M19").Value = 0.02222222222222222
$ws.Range("N19").Value = 0.0007936507936507937
$ws.Range("O19").Value = 0.06349206349206349
$ws.Range("S19").Value = 0.09603174603174604
